# Update LR-pairs data with refreshed TPM output.
# Final state: only 3 data rows remain (rows 2-4); sending/ligand/receptor
# columns are FAPs/Tnfsf11/Tnfrsf11b in every row, and the target-cluster
# column cycles through ECs, FAPs, MuSCs with refreshed numeric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused trailing rows (old rows 5, 6, 7) so only
# header + 3 data rows remain. This also removes the last references to
# "Resolving-Mac" (previously only used in those rows).
$ws.Rows.Item(5).Resize(3).Delete() | Out-Null

# Clear out the remaining references to "MuSCs" (old rows 3 & 4, column A/D)
# before re-introducing it later, further down the sheet, so the shared
# string table regenerates it in the right relative position.
$ws.Range("D3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"

# --- Row 2: FAPs -> Tnfsf11/Tnfrsf11b -> ECs ---
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf11"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.690534333333333
$ws.Range("H2").Value = 5.071603
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06861733333333334
$ws.Range("N2").Value = 0.205852
$ws.Range("O2").Value = 0.01654048691795588
$ws.Range("P2").Value = 0.01654048691795588
$ws.Range("Q2").Value = 0.1159999578617778
$ws.Range("R2").Value = 1.043999620756
$ws.Range("S2").Value = 0.01654048691795588
$ws.Range("T2").Value = 0.01654048691795588

# --- Row 3: FAPs -> Tnfsf11/Tnfrsf11b -> FAPs ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf11"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.690534333333333
$ws.Range("H3").Value = 5.071603
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.776574666666666
$ws.Range("N3").Value = 11.329724
$ws.Range("O3").Value = 0.9103586635352137
$ws.Range("P3").Value = 0.9103586635352137
$ws.Range("Q3").Value = 6.384429136396887
$ws.Range("R3").Value = 57.45986222757199
$ws.Range("S3").Value = 0.9103586635352137
$ws.Range("T3").Value = 0.9103586635352137

# --- Row 4: FAPs -> Tnfsf11/Tnfrsf11b -> MuSCs ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf11"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.690534333333333
$ws.Range("H4").Value = 5.071603
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.303255
$ws.Range("N4").Value = 0.909765
$ws.Range("O4").Value = 0.07310084954683041
$ws.Range("P4").Value = 0.07310084954683042
$ws.Range("Q4").Value = 0.512662989255
$ws.Range("R4").Value = 4.613966903295
$ws.Range("S4").Value = 0.07310084954683041
$ws.Range("T4").Value = 0.07310084954683042
